$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.62
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.49

# Row 3
$ws.Range("G3").Value = 1.67
$ws.Range("Q3").Value = 2.15
$ws.Range("R3").Value = 1.61

# Row 4
$ws.Range("K4").Value = 1.92
$ws.Range("L4").Value = 2.87
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.47

# Row 5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5

# Row 6
$ws.Range("G6").Value = 5.5
$ws.Range("I6").Value = 1.73
$ws.Range("J6").Value = 6
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.63
$ws.Range("R6").Value = 1.5
$ws.Range("W6").Value = 11
$ws.Range("X6").Value = 26
$ws.Range("Z6").Value = 67
$ws.Range("AH6").Value = 5
$ws.Range("AN6").Value = 7
$ws.Range("AO6").Value = 34
$ws.Range("AR6").Value = 201
$ws.Range("AU6").Value = 10
$ws.Range("AX6").Value = 9.5

# Row 7
$ws.Range("G7").Value = 1.83

# Row 11
$ws.Range("U11").Value = 2.37
$ws.Range("V11").Value = 1.5

# Row 12
$ws.Range("G12").Value = 1.65
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 2.38
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("U12").Value = 2.37
$ws.Range("V12").Value = 1.5
$ws.Range("W12").Value = 5
$ws.Range("X12").Value = 6.5
$ws.Range("Z12").Value = 12
$ws.Range("AH12").Value = 12
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 21
$ws.Range("AM12").Value = 67
$ws.Range("AN12").Value = 3.4
$ws.Range("AO12").Value = 9
$ws.Range("AU12").Value = 10
$ws.Range("AW12").Value = 7
$ws.Range("AZ12").Value = 151
$ws.Range("BA12").Value = 201

# Row 13
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 1.98
$ws.Range("R13").Value = 1.88
$ws.Range("U13").Value = 1.69
